$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AX = column 50. Match width of neighboring column AW (49).
$ws.Columns.Item(50).ColumnWidth = $ws.Columns.Item(49).ColumnWidth

# Row 1 header: literal text date (not an Excel date value), same style as rest of header row.
$ws.Cells.Item(1, 50).Value = "'2024/10/28"
$ws.Cells.Item(1, 50).Font.Name = $ws.Cells.Item(1, 2).Font.Name

# Data rows 2-53: numeric values, font + fill matched to the style used in the diff
# style 1 = plain (font only), style 2 = yellow fill, style 3 = light-blue fill
$yellow = 65535
$blue = 15128749

$ws.Cells.Item(2, 50).Value = 188
$ws.Cells.Item(2, 50).Font.Name = $ws.Cells.Item(2, 2).Font.Name
$ws.Cells.Item(3, 50).Value = 223.3
$ws.Cells.Item(3, 50).Font.Name = $ws.Cells.Item(3, 2).Font.Name
$ws.Cells.Item(4, 50).Value = 264.5
$ws.Cells.Item(4, 50).Font.Name = $ws.Cells.Item(4, 2).Font.Name
$ws.Cells.Item(5, 50).Value = 112
$ws.Cells.Item(5, 50).Font.Name = $ws.Cells.Item(5, 2).Font.Name
$ws.Cells.Item(5, 50).Interior.Color = $yellow
$ws.Cells.Item(6, 50).Value = 227.4
$ws.Cells.Item(6, 50).Font.Name = $ws.Cells.Item(6, 2).Font.Name
$ws.Cells.Item(7, 50).Value = 140.9
$ws.Cells.Item(7, 50).Font.Name = $ws.Cells.Item(7, 2).Font.Name
$ws.Cells.Item(8, 50).Value = 268.8
$ws.Cells.Item(8, 50).Font.Name = $ws.Cells.Item(8, 2).Font.Name
$ws.Cells.Item(9, 50).Value = 118.4
$ws.Cells.Item(9, 50).Font.Name = $ws.Cells.Item(9, 2).Font.Name
$ws.Cells.Item(9, 50).Interior.Color = $yellow
$ws.Cells.Item(10, 50).Value = 124.6
$ws.Cells.Item(10, 50).Font.Name = $ws.Cells.Item(10, 2).Font.Name
$ws.Cells.Item(10, 50).Interior.Color = $yellow
$ws.Cells.Item(11, 50).Value = 136.6
$ws.Cells.Item(11, 50).Font.Name = $ws.Cells.Item(11, 2).Font.Name
$ws.Cells.Item(11, 50).Interior.Color = $blue
$ws.Cells.Item(12, 50).Value = 168.8
$ws.Cells.Item(12, 50).Font.Name = $ws.Cells.Item(12, 2).Font.Name
$ws.Cells.Item(13, 50).Value = 155.7
$ws.Cells.Item(13, 50).Font.Name = $ws.Cells.Item(13, 2).Font.Name
$ws.Cells.Item(14, 50).Value = 123.3
$ws.Cells.Item(14, 50).Font.Name = $ws.Cells.Item(14, 2).Font.Name
$ws.Cells.Item(14, 50).Interior.Color = $yellow
$ws.Cells.Item(15, 50).Value = 196
$ws.Cells.Item(15, 50).Font.Name = $ws.Cells.Item(15, 2).Font.Name
$ws.Cells.Item(16, 50).Value = 166.5
$ws.Cells.Item(16, 50).Font.Name = $ws.Cells.Item(16, 2).Font.Name
$ws.Cells.Item(17, 50).Value = 184.8
$ws.Cells.Item(17, 50).Font.Name = $ws.Cells.Item(17, 2).Font.Name
$ws.Cells.Item(18, 50).Value = 147.8
$ws.Cells.Item(18, 50).Font.Name = $ws.Cells.Item(18, 2).Font.Name
$ws.Cells.Item(19, 50).Value = 148.5
$ws.Cells.Item(19, 50).Font.Name = $ws.Cells.Item(19, 2).Font.Name
$ws.Cells.Item(20, 50).Value = 160.4
$ws.Cells.Item(20, 50).Font.Name = $ws.Cells.Item(20, 2).Font.Name
$ws.Cells.Item(21, 50).Value = 162.2
$ws.Cells.Item(21, 50).Font.Name = $ws.Cells.Item(21, 2).Font.Name
$ws.Cells.Item(22, 50).Value = 182.9
$ws.Cells.Item(22, 50).Font.Name = $ws.Cells.Item(22, 2).Font.Name
$ws.Cells.Item(23, 50).Value = 153.3
$ws.Cells.Item(23, 50).Font.Name = $ws.Cells.Item(23, 2).Font.Name
$ws.Cells.Item(24, 50).Value = 123.8
$ws.Cells.Item(24, 50).Font.Name = $ws.Cells.Item(24, 2).Font.Name
$ws.Cells.Item(24, 50).Interior.Color = $yellow
$ws.Cells.Item(25, 50).Value = 117.2
$ws.Cells.Item(25, 50).Font.Name = $ws.Cells.Item(25, 2).Font.Name
$ws.Cells.Item(25, 50).Interior.Color = $yellow
$ws.Cells.Item(26, 50).Value = 134.7
$ws.Cells.Item(26, 50).Font.Name = $ws.Cells.Item(26, 2).Font.Name
$ws.Cells.Item(26, 50).Interior.Color = $blue
$ws.Cells.Item(27, 50).Value = 136.4
$ws.Cells.Item(27, 50).Font.Name = $ws.Cells.Item(27, 2).Font.Name
$ws.Cells.Item(27, 50).Interior.Color = $blue
$ws.Cells.Item(28, 50).Value = 197.6
$ws.Cells.Item(28, 50).Font.Name = $ws.Cells.Item(28, 2).Font.Name
$ws.Cells.Item(29, 50).Value = 122.2
$ws.Cells.Item(29, 50).Font.Name = $ws.Cells.Item(29, 2).Font.Name
$ws.Cells.Item(29, 50).Interior.Color = $yellow
$ws.Cells.Item(30, 50).Value = 158.8
$ws.Cells.Item(30, 50).Font.Name = $ws.Cells.Item(30, 2).Font.Name
$ws.Cells.Item(31, 50).Value = 184.9
$ws.Cells.Item(31, 50).Font.Name = $ws.Cells.Item(31, 2).Font.Name
$ws.Cells.Item(32, 50).Value = 166.5
$ws.Cells.Item(32, 50).Font.Name = $ws.Cells.Item(32, 2).Font.Name
$ws.Cells.Item(33, 50).Value = 125.9
$ws.Cells.Item(33, 50).Font.Name = $ws.Cells.Item(33, 2).Font.Name
$ws.Cells.Item(33, 50).Interior.Color = $blue
$ws.Cells.Item(34, 50).Value = 125.9
$ws.Cells.Item(34, 50).Font.Name = $ws.Cells.Item(34, 2).Font.Name
$ws.Cells.Item(34, 50).Interior.Color = $blue
$ws.Cells.Item(35, 50).Value = 209
$ws.Cells.Item(35, 50).Font.Name = $ws.Cells.Item(35, 2).Font.Name
$ws.Cells.Item(36, 50).Value = 177.1
$ws.Cells.Item(36, 50).Font.Name = $ws.Cells.Item(36, 2).Font.Name
$ws.Cells.Item(37, 50).Value = 142.4
$ws.Cells.Item(37, 50).Font.Name = $ws.Cells.Item(37, 2).Font.Name
$ws.Cells.Item(38, 50).Value = 134.2
$ws.Cells.Item(38, 50).Font.Name = $ws.Cells.Item(38, 2).Font.Name
$ws.Cells.Item(38, 50).Interior.Color = $blue
$ws.Cells.Item(39, 50).Value = 300.2
$ws.Cells.Item(39, 50).Font.Name = $ws.Cells.Item(39, 2).Font.Name
$ws.Cells.Item(40, 50).Value = 157.5
$ws.Cells.Item(40, 50).Font.Name = $ws.Cells.Item(40, 2).Font.Name
$ws.Cells.Item(41, 50).Value = 148.7
$ws.Cells.Item(41, 50).Font.Name = $ws.Cells.Item(41, 2).Font.Name
$ws.Cells.Item(42, 50).Value = 124.2
$ws.Cells.Item(42, 50).Font.Name = $ws.Cells.Item(42, 2).Font.Name
$ws.Cells.Item(42, 50).Interior.Color = $yellow
$ws.Cells.Item(43, 50).Value = 184.9
$ws.Cells.Item(43, 50).Font.Name = $ws.Cells.Item(43, 2).Font.Name
$ws.Cells.Item(44, 50).Value = 129.7
$ws.Cells.Item(44, 50).Font.Name = $ws.Cells.Item(44, 2).Font.Name
$ws.Cells.Item(44, 50).Interior.Color = $blue
$ws.Cells.Item(45, 50).Value = 148.5
$ws.Cells.Item(45, 50).Font.Name = $ws.Cells.Item(45, 2).Font.Name
$ws.Cells.Item(46, 50).Value = 127.5
$ws.Cells.Item(46, 50).Font.Name = $ws.Cells.Item(46, 2).Font.Name
$ws.Cells.Item(46, 50).Interior.Color = $blue
$ws.Cells.Item(47, 50).Value = 127.9
$ws.Cells.Item(47, 50).Font.Name = $ws.Cells.Item(47, 2).Font.Name
$ws.Cells.Item(47, 50).Interior.Color = $blue
$ws.Cells.Item(48, 50).Value = 177.8
$ws.Cells.Item(48, 50).Font.Name = $ws.Cells.Item(48, 2).Font.Name
$ws.Cells.Item(49, 50).Value = 214.1
$ws.Cells.Item(49, 50).Font.Name = $ws.Cells.Item(49, 2).Font.Name
$ws.Cells.Item(50, 50).Value = 150.2
$ws.Cells.Item(50, 50).Font.Name = $ws.Cells.Item(50, 2).Font.Name
$ws.Cells.Item(51, 50).Value = 147
$ws.Cells.Item(51, 50).Font.Name = $ws.Cells.Item(51, 2).Font.Name
$ws.Cells.Item(52, 50).Value = 168.2
$ws.Cells.Item(52, 50).Font.Name = $ws.Cells.Item(52, 2).Font.Name
$ws.Cells.Item(53, 50).Value = 201.7
$ws.Cells.Item(53, 50).Font.Name = $ws.Cells.Item(53, 2).Font.Name
